$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update price (D) and volume (E) columns for changed rows.
# D-column prices are forced to text via a leading apostrophe (quote prefix)
# so Excel keeps them as strings instead of auto-converting to numbers,
# then the style is reset to Normal so no stray number-format style is left behind.
$ws.Range("D2").Value = "'40.110.82"
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = "  +0.32%  "
$ws.Range("D3").Value = "'2.225.50"
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = "  -0.56%  "
$ws.Range("D5").Value = "'290.87"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  -0.60%  "
$ws.Range("D6").Value = "'87.58"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  +1.35%  "
$ws.Range("E7").Value = "  -0.41%  "
$ws.Range("E8").Value = "  +0.02%  "
$ws.Range("E9").Value = "  +0.06%  "
$ws.Range("D10").Value = "'30.39"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  -1.89%  "
$ws.Range("E11").Value = "  -2.22%  "
$ws.Range("E12").Value = "  +3.00%  "
$ws.Range("D13").Value = "'6.48"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "  +0.79%  "
$ws.Range("D14").Value = "'2.570.85"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  -0.52%  "
$ws.Range("D15").Value = "'13.91"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "  -2.01%  "
$ws.Range("D16").Value = "'2.223.39"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "  -0.84%  "
$ws.Range("D17").Value = "'0.728"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "  -0.40%  "
$ws.Range("D18").Value = "'40.046.27"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "  +0.36%  "
$ws.Range("D19").Value = "'0.0₃0887"
$ws.Range("D19").Style = "Normal"
$ws.Range("D20").Value = "'11.37"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  +7.71%  "
$ws.Range("E21").Value = "  -0.26%  "
$ws.Range("D22").Value = "'65.70"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  +0.03%  "
$ws.Range("D23").Value = "'237.01"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  +0.50%  "
$ws.Range("E24").Value = "  -0.04%  "
$ws.Range("D25").Value = "'2.45"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  +0.21%  "
$ws.Range("E26").Value = "  -1.27%  "
$ws.Range("D27").Value = "'22.72"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  -1.41%  "
$ws.Range("E28").Value = "  -1.67%  "
$ws.Range("D29").Value = "'9.23"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  -0.50%  "
$ws.Range("D30").Value = "'155.96"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  +1.52%  "
$ws.Range("E31").Value = "  -6.48%  "
$ws.Range("E32").Value = "  -0.15%  "
$ws.Range("D33").Value = "'4.95"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  +1.37%  "
$ws.Range("D34").Value = "'0.0721"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  +0.96%  "
$ws.Range("E35").Value = "  +7.34%  "
$ws.Range("D36").Value = "'2.35"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  -1.32%  "
$ws.Range("E37").Value = "  +0.57%  "
$ws.Range("D39").Value = "'0.0982"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  -2.05%  "
$ws.Range("E40").Value = "  +1.36%  "
$ws.Range("D41").Value = "'2.130.87"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  +7.96%  "
$ws.Range("D42").Value = "'3.87"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  +2.33%  "
$ws.Range("D43").Value = "'18.66"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  +14.92%  "
$ws.Range("E44").Value = "  -3.97%  "
$ws.Range("E45").Value = "  -1.15%  "
$ws.Range("D46").Value = "'9.86"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  +0.55%  "
$ws.Range("D47").Value = "'2.66"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  +3.61%  "
$ws.Range("D48").Value = "'2.437.23"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  -0.74%  "
$ws.Range("E49").Value = "  +0.06%  "

# Rows 50 and 51 swapped rank (BitcoinSV now ranks above TrustWalletToken) with updated values
$ws.Range("B50").Value = "BitcoinSV"
$ws.Range("C50").Value = "https://coinranking.com/coin/VcMY11NONHSA0+bitcoinsv-bsv"
$ws.Range("D50").Value = "'69.34"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  -2.48%  "

$ws.Range("B51").Value = "TrustWalletToken"
$ws.Range("C51").Value = "https://coinranking.com/coin/Hm3OlynlC+trustwallettoken-twt"
$ws.Range("D51").Value = "'1.10"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  +2.05%  "
